$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab title from "Through 2022-10-01" to "Through 2022-10-02"
$ws.Name = "Through 2022-10-02"

# Update header label in I1 (shared string) from "2022 (through 10-01)" to "2022 (through 10-02)"
$ws.Range("I1").Value = "2022 (through 10-02)"

# Update October value for 2022 column (I11): 5 -> 9
$ws.Range("I11").Value = 9

# Update Total value for 2022 column (I14): 1287 -> 1291
$ws.Range("I14").Value = 1291
